$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "26.775.08"
Set-TextValue $ws.Cells.Item(2, 5) "  +1.47%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "1.721.27"
Set-TextValue $ws.Cells.Item(3, 5) "  +0.19%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "1.001"
Set-TextValue $ws.Cells.Item(4, 5) "  +0.34%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "239.92"
Set-TextValue $ws.Cells.Item(5, 5) "  -0.72%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 5) "  +0.31%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.4742"
Set-TextValue $ws.Cells.Item(7, 5) "  -2.48%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.2556"
Set-TextValue $ws.Cells.Item(8, 5) "  -1.06%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.06119"
Set-TextValue $ws.Cells.Item(9, 5) "  -0.72%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "1.723.25"
Set-TextValue $ws.Cells.Item(10, 5) "  +0.01%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "15.84"
Set-TextValue $ws.Cells.Item(11, 5) "  +2.32%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.06876"
Set-TextValue $ws.Cells.Item(12, 5) "  -1.00%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.5965"
Set-TextValue $ws.Cells.Item(13, 5) "  -0.13%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "4.404"
Set-TextValue $ws.Cells.Item(14, 5) "  -1.92%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "76.40"
Set-TextValue $ws.Cells.Item(15, 5) "  -0.27%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 5) "  +0.31%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "26.612.26"
Set-TextValue $ws.Cells.Item(17, 5) "  +0.84%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "1.002"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.39%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "0.000007058"
Set-TextValue $ws.Cells.Item(19, 5) "  -0.60%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "11.23"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.08%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "1.945.07"
Set-TextValue $ws.Cells.Item(21, 5) "  -0.34%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "4.368"
Set-TextValue $ws.Cells.Item(22, 5) "  -1.06%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "8.309"
Set-TextValue $ws.Cells.Item(23, 5) "  -1.65%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "5.015"
Set-TextValue $ws.Cells.Item(24, 5) "  -1.07%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "140.22"
Set-TextValue $ws.Cells.Item(25, 5) "  +2.41%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "15.12"
Set-TextValue $ws.Cells.Item(26, 5) "  -0.57%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "1.772"
Set-TextValue $ws.Cells.Item(27, 5) "  +2.14%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "106.10"
Set-TextValue $ws.Cells.Item(28, 5) "  +0.18%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 5) "  -2.42%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "3.903"
Set-TextValue $ws.Cells.Item(30, 5) "  +0.56%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "0.07856"
Set-TextValue $ws.Cells.Item(31, 5) "  -1.28%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "3.616"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.19%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "0.04527"
Set-TextValue $ws.Cells.Item(33, 5) "  +1.87%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 5) "  -0.04%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "0.9895"
Set-TextValue $ws.Cells.Item(35, 5) "  -0.55%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "0.6095"
Set-TextValue $ws.Cells.Item(36, 5) "  -1.52%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "0.9144"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.95%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "2.484"
Set-TextValue $ws.Cells.Item(38, 5) "  +4.57%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "1.960"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.40%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "1.000"
Set-TextValue $ws.Cells.Item(40, 5) "  +0.29%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "5.715"
Set-TextValue $ws.Cells.Item(41, 5) "  +4.76%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "0.01476"
Set-TextValue $ws.Cells.Item(42, 5) "  +0.00%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "99.87"
Set-TextValue $ws.Cells.Item(43, 5) "  +0.74%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.3780"
Set-TextValue $ws.Cells.Item(44, 5) "  -0.94%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "6.686"
Set-TextValue $ws.Cells.Item(45, 5) "  -2.09%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "0.1137"
Set-TextValue $ws.Cells.Item(46, 5) "  -1.18%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "0.05352"
Set-TextValue $ws.Cells.Item(47, 5) "  -0.09%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "7.816"
Set-TextValue $ws.Cells.Item(48, 5) "  +1.38%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "29.62"
Set-TextValue $ws.Cells.Item(49, 5) "  -2.52%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "1.223"
Set-TextValue $ws.Cells.Item(50, 5) "  +0.49%  "

# Row 51 (Aave -> TrueUSD)
Set-TextValue $ws.Cells.Item(51, 2) "TrueUSD"
Set-TextValue $ws.Cells.Item(51, 3) "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
Set-TextValue $ws.Cells.Item(51, 4) "1.004"
Set-TextValue $ws.Cells.Item(51, 5) "  +0.35%  "